# "Generate Report for handback"
#
# The 306d83a0-e56d-4b01-95b8-a32b948f5fcc.md file (rows 3/4 on each locale
# sheet) has now been handed back, so:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The per-locale sheets grow a "Latest Target File" / "Latest Handback
#     File" pair of hyperlinks (columns E/F) mirroring the ones already
#     present for the first file in row 2
#   - "Latest Handback DateTime" (col G) is stamped with the handback time
#   - "Handoff Reason" (col H) flips from "Ignored" to "Include"
#
# This mirrors what Excel does when a user fills the new handback info in
# for row 3 / row 4 across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the Status column (B) and the per locale columns
# mirror it (C here == zh-cn status, since "de-de" column value equals
# the same status text too).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack
$wsOverview.Range("B4").Value = $statusHandedBack
$wsOverview.Range("C4").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B3").Value = $statusHandedBack
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e5e18c77ce81350ca97d6516f7174c905528c86d/e2e/306d83a0-e56d-4b01-95b8-a32b948f5fcc.md", "", "", "306d83a0-e56d-4b01-95b8-a32b948f5fcc.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4a24269b214b59a58208483ed8d937fbf7e0af57/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/qimu/306d83a0-e56d-4b01-95b8-a32b948f5fcc.dae33526312388e6fe346ae598396f24799588e3.zh-cn.xlf", "", "", "306d83a0-e56d-4b01-95b8-a32b948f5fcc.dae33526312388e6fe346ae598396f24799588e3.zh-cn.xlf")
$wsZh.Range("G3").Value = "2016-01-25 03:42:47"
$wsZh.Range("H3").Value = "Include"

$wsZh.Range("B4").Value = $statusHandedBack
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e5e18c77ce81350ca97d6516f7174c905528c86d/e2e/306d83a0-e56d-4b01-95b8-a32b948f5fcc.md", "", "", "306d83a0-e56d-4b01-95b8-a32b948f5fcc.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4a24269b214b59a58208483ed8d937fbf7e0af57/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/qimu/306d83a0-e56d-4b01-95b8-a32b948f5fcc.dae33526312388e6fe346ae598396f24799588e3.zh-cn.xlf", "", "", "306d83a0-e56d-4b01-95b8-a32b948f5fcc.dae33526312388e6fe346ae598396f24799588e3.zh-cn.xlf")
$wsZh.Range("G4").Value = "2016-01-25 03:42:47"
$wsZh.Range("H4").Value = "Include"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B3").Value = $statusHandedBack
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e5e18c77ce81350ca97d6516f7174c905528c86d/e2e/306d83a0-e56d-4b01-95b8-a32b948f5fcc.md", "", "", "306d83a0-e56d-4b01-95b8-a32b948f5fcc.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2107ee31001ce2e6b20019af92194920fd09b5b0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/qimu/306d83a0-e56d-4b01-95b8-a32b948f5fcc.dae33526312388e6fe346ae598396f24799588e3.de-de.xlf", "", "", "306d83a0-e56d-4b01-95b8-a32b948f5fcc.dae33526312388e6fe346ae598396f24799588e3.de-de.xlf")
$wsDe.Range("G3").Value = "2016-01-25 03:43:04"
$wsDe.Range("H3").Value = "Include"

$wsDe.Range("B4").Value = $statusHandedBack
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e5e18c77ce81350ca97d6516f7174c905528c86d/e2e/306d83a0-e56d-4b01-95b8-a32b948f5fcc.md", "", "", "306d83a0-e56d-4b01-95b8-a32b948f5fcc.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2107ee31001ce2e6b20019af92194920fd09b5b0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/qimu/306d83a0-e56d-4b01-95b8-a32b948f5fcc.dae33526312388e6fe346ae598396f24799588e3.de-de.xlf", "", "", "306d83a0-e56d-4b01-95b8-a32b948f5fcc.dae33526312388e6fe346ae598396f24799588e3.de-de.xlf")
$wsDe.Range("G4").Value = "2016-01-25 03:43:04"
$wsDe.Range("H4").Value = "Include"
